# Scheduled Universalis market-data refresh for the Odin_Profits leve-profit tracker.
# For every touched leve row this re-prices columns H:N:
#   H  currentAveragePrice      I  currentAveragePriceNQ   J  currentAveragePriceHQ
#   K  LevePriceNQ              L  LevePriceHQ
#   M  LeveProfitNQ             N  LeveProfitHQ
$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")

# Row 33: "Glazed and Confused" (Clear Glass Lens)
$ws.Range("H33").Value = 748.13635
$ws.Range("I33").Value = 737.93335
$ws.Range("J33").Value = 770
$ws.Range("K33").Value = 737.93335
$ws.Range("L33").Value = 770
$ws.Range("M33").Value = -508.93335
$ws.Range("N33").Value = -1228

# Row 43: "Growing Is Knowing" (Growth Formula Gamma)
$ws.Range("H43").Value = 1761.2858
$ws.Range("I43").Value = 880
$ws.Range("J43").Value = 1908.1666
$ws.Range("K43").Value = 880
$ws.Range("L43").Value = 1908.1666
$ws.Range("M43").Value = -811
$ws.Range("N43").Value = -2046.1666

# Row 74: "Adhesive of Antipathy" (Wing Glue)
$ws.Range("H74").Value = 11594.782
$ws.Range("I74").Value = 4075.2
$ws.Range("J74").Value = 13683.556
$ws.Range("K74").Value = 4075.2
$ws.Range("L74").Value = 13683.556
$ws.Range("M74").Value = -3139.2
$ws.Range("N74").Value = -15555.556

# Row 77: "It's Gonna Grow Back (L)" (Wing Glue)
$ws.Range("H77").Value = 11594.782
$ws.Range("I77").Value = 4075.2
$ws.Range("J77").Value = 13683.556
$ws.Range("K77").Value = 20376
$ws.Range("L77").Value = 68417.78
$ws.Range("M77").Value = -15696
$ws.Range("N77").Value = -77777.78

# Row 135: "For Tired Minds" (Grade 1 Gemsap of Intelligence)
$ws.Range("H135").Value = 11915.083
$ws.Range("I135").Value = 5996.6
$ws.Range("K135").Value = 53969.4
$ws.Range("M135").Value = -51434.4

# Row 136: "I Like Big Brush and I Cannot Lie" (Dark Mahogany Round Brush)
$ws.Range("H136").Value = 148986
$ws.Range("J136").Value = 148986
$ws.Range("L136").Value = 148986
$ws.Range("N136").Value = -159186

# Row 137: "Cutting Edge of Culinary Quality" (Magnesia Whetstone)
$ws.Range("H137").Value = 12453.412
$ws.Range("J137").Value = 15693.625
$ws.Range("L137").Value = 47080.875
$ws.Range("N137").Value = -52180.875

# Row 138: "All-night Crafting" (Cunning Craftsman's Tisane)
$ws.Range("H138").Value = 5156.7827
$ws.Range("I138").Value = 933.6
$ws.Range("J138").Value = 6329.8887
$ws.Range("K138").Value = 2800.8
$ws.Range("L138").Value = 18989.6661
$ws.Range("M138").Value = 2339.2
$ws.Range("N138").Value = -29269.6661

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")

# Row 32: "Ingot We Trust" (Steel Ingot)
$ws.Range("H32").Value = 2608.2354
$ws.Range("I32").Value = 905.5454999999999
$ws.Range("J32").Value = 5729.8335
$ws.Range("K32").Value = 905.5454999999999
$ws.Range("L32").Value = 5729.8335
$ws.Range("M32").Value = -618.5454999999999
$ws.Range("N32").Value = -6303.8335

# Row 43: "They've Got Legs" (Steel Sabatons)
$ws.Range("H43").Value = 28727.625
$ws.Range("J43").Value = 26745.857
$ws.Range("L43").Value = 26745.857
$ws.Range("N43").Value = -27371.857

# Row 45: "Hollow Hallmarks" (Mythril Ingot)
$ws.Range("H45").Value = 2911.3684
$ws.Range("I45").Value = 2929
$ws.Range("K45").Value = 2929
$ws.Range("M45").Value = -2552

# Row 61: "Dealing with the Tough Stuff" (Cobalt Ingot)
$ws.Range("H61").Value = 6755
$ws.Range("I61").Value = 11141.385
$ws.Range("J61").Value = 2681.9285
$ws.Range("K61").Value = 11141.385
$ws.Range("L61").Value = 2681.9285
$ws.Range("M61").Value = -10929.385
$ws.Range("N61").Value = -3105.9285

# Row 74: "As the Bolt Flies" (Titanium Nugget)
$ws.Range("H74").Value = 6991.1113
$ws.Range("I74").Value = 10240.667
$ws.Range("J74").Value = 5366.3335
$ws.Range("K74").Value = 10240.667
$ws.Range("L74").Value = 5366.3335
$ws.Range("M74").Value = -9366.666999999999
$ws.Range("N74").Value = -7114.3335

# Row 77: "Heavy Metal Banned (L)" (Titanium Nugget)
$ws.Range("H77").Value = 6991.1113
$ws.Range("I77").Value = 10240.667
$ws.Range("J77").Value = 5366.3335
$ws.Range("K77").Value = 51203.335
$ws.Range("L77").Value = 26831.6675
$ws.Range("M77").Value = -46835.335
$ws.Range("N77").Value = -35567.6675

# Row 136: "Metal with Mettle" (Cobalt Tungsten Ingot)
$ws.Range("H136").Value = 6755
$ws.Range("I136").Value = 11141.385
$ws.Range("J136").Value = 2681.9285
$ws.Range("K136").Value = 33424.155
$ws.Range("L136").Value = 8045.7855
$ws.Range("M136").Value = -30874.155
$ws.Range("N136").Value = -13145.7855

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")

# Row 99: "Meddle in Metal" (Oroshigane Ingot)
$ws.Range("H99").Value = 7480.1113
$ws.Range("I99").Value = 4206.25
$ws.Range("J99").Value = 10099.2
$ws.Range("K99").Value = 4206.25
$ws.Range("L99").Value = 10099.2
$ws.Range("M99").Value = -2708.25
$ws.Range("N99").Value = -13095.2

# Row 134: "Ruthenium Supremium" (Ruthenium Ingot)
$ws.Range("H134").Value = 8284.522999999999
$ws.Range("I134").Value = 7211.25
$ws.Range("K134").Value = 21633.75
$ws.Range("M134").Value = -19098.75

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")

# Row 94: "Beech, Please" (Beech Lumber)
$ws.Range("H94").Value = 3153.24
$ws.Range("I94").Value = 1280.75
$ws.Range("J94").Value = 6482.1113
$ws.Range("K94").Value = 1280.75
$ws.Range("L94").Value = 6482.1113
$ws.Range("M94").Value = -829.75
$ws.Range("N94").Value = -7384.1113

# Row 134: "Wood You Be Quiet" (Ceiba Lumber)
$ws.Range("H134").Value = 11637.482
$ws.Range("I134").Value = 10432
$ws.Range("J134").Value = 12096.714
$ws.Range("K134").Value = 31296
$ws.Range("L134").Value = 36290.142
$ws.Range("M134").Value = -28761
$ws.Range("N134").Value = -41360.142

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")

# Row 4: "In Hot Water" (Boiled Egg)
$ws.Range("H4").Value = 6548212.5
$ws.Range("I4").Value = 7857544.5
$ws.Range("K4").Value = 23572633.5
$ws.Range("M4").Value = -23572521.5

# Row 6: "Meat-lover's Special" (Marmot Steak)
$ws.Range("H6").Value = 814
$ws.Range("I6").Value = 682.2222
$ws.Range("J6").Value = 2000
$ws.Range("K6").Value = 2046.6666
$ws.Range("L6").Value = 6000
$ws.Range("M6").Value = -1933.6666
$ws.Range("N6").Value = -6226

# Row 7: "It's Always Sunny in Vylbrand" (Raisins)
$ws.Range("H7").Value = 444.75
$ws.Range("J7").Value = 149
$ws.Range("L7").Value = 447
$ws.Range("N7").Value = -671

# Row 9: "Jack of All Plates" (Jack-o'-lantern)
$ws.Range("H9").Value = 167074.73
$ws.Range("I9").Value = 3249
$ws.Range("J9").Value = 203480.44
$ws.Range("K9").Value = 9747
$ws.Range("L9").Value = 610441.3200000001
$ws.Range("M9").Value = -9523
$ws.Range("N9").Value = -610889.3200000001

# Row 10: "A Real Fungi" (Chanterelle Saute)
$ws.Range("H10").Value = 116.85714
$ws.Range("I10").Value = 95.59999999999999
$ws.Range("K10").Value = 286.8
$ws.Range("M10").Value = -147.8

# Row 11: "Putting the Squeeze On" (Orange Juice)
$ws.Range("H11").Value = 17.5
$ws.Range("I11").Value = 17.5
$ws.Range("K11").Value = 52.5
$ws.Range("M11").Value = 87.5

# Row 16: "Go Ahead and Dig In" (Mole Loaf)
$ws.Range("H16").Value = 2000
$ws.Range("I16").Value = 2000
$ws.Range("K16").Value = 6000
$ws.Range("M16").Value = -5827

# Row 24: "Rustic Repast" (Chicken and Mushrooms)
$ws.Range("H24").Value = 3999.6667

# Row 29: "For Crumbs' Sake" (Honey Muffin)
$ws.Range("H29").Value = 132.83333
$ws.Range("I29").Value = 132.83333
$ws.Range("K29").Value = 398.49999
$ws.Range("M29").Value = -121.49999

# Row 31: "Food Fight" (Shepherd's Pie)
$ws.Range("H31").Value = 2000
$ws.Range("I31").Value = 2000
$ws.Range("K31").Value = 6000
$ws.Range("M31").Value = -5712

# Row 98: "Sweet Kiss of Death" (Rice Vinegar)
$ws.Range("H98").Value = 1756.8572
$ws.Range("I98").Value = 3000.5
$ws.Range("J98").Value = 1259.4
$ws.Range("K98").Value = 9001.5
$ws.Range("L98").Value = 3778.2
$ws.Range("M98").Value = -7503.5
$ws.Range("N98").Value = -6774.200000000001

# Row 126: "Imperial Palate" (Glory Be Soup)
$ws.Range("H126").Value = 6677.2
$ws.Range("I126").Value = 6426
$ws.Range("J126").Value = 7054
$ws.Range("K126").Value = 19278
$ws.Range("L126").Value = 21162
$ws.Range("M126").Value = -14338
$ws.Range("N126").Value = -31042

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")

# Row 124: "The Sage's Successor" (Pewter Pendulums)
$ws.Range("H124").Value = 62856.285
$ws.Range("J124").Value = 62856.285
$ws.Range("L124").Value = 62856.285
$ws.Range("N124").Value = -72676.285

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")

# Row 30: "Packing a Punch" (Goatskin Cesti)
$ws.Range("H30").Value = 857
$ws.Range("I30").Value = 857
$ws.Range("J30").Value = 0
$ws.Range("K30").Value = 857
$ws.Range("L30").Value = 0
$ws.Range("M30").Value = -749
$ws.Range("N30").ClearContents()

# Row 46: "Supply Side Logic" (Boar Leather)
$ws.Range("H46").Value = 843.1667
$ws.Range("I46").Value = 726.6667
$ws.Range("J46").Value = 959.6667
$ws.Range("K46").Value = 726.6667
$ws.Range("L46").Value = 959.6667
$ws.Range("M46").Value = -538.6667
$ws.Range("N46").Value = -1335.6667

# Row 61: "Spelling Me Softly" (Raptor Leather)
$ws.Range("H61").Value = 9787.615
$ws.Range("I61").Value = 8124.0835
$ws.Range("K61").Value = 8124.0835
$ws.Range("M61").Value = -7922.0835

# Row 113: "Peace in Rest" (Atrociraptor Leather)
$ws.Range("H113").Value = 9787.615
$ws.Range("I113").Value = 8124.0835
$ws.Range("K113").Value = 8124.0835
$ws.Range("M113").Value = -5954.0835

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")

# Row 126: "A Polished Purchase" (Snow Linen)
$ws.Range("H126").Value = 4949.2354
$ws.Range("I126").Value = 3521.7778
$ws.Range("J126").Value = 6555.125
$ws.Range("K126").Value = 10565.3334
$ws.Range("L126").Value = 19665.375
$ws.Range("M126").Value = -8095.3334
$ws.Range("N126").Value = -24605.375

# Row 136: "Weaving the Envelope" (Sarcenet Cloth)
$ws.Range("H136").Value = 68533.164
$ws.Range("I136").Value = 103155.3
$ws.Range("K136").Value = 309465.9
$ws.Range("M136").Value = -306915.9
